$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 17-21 (rows removed in this edit)
$ws.Rows("17:21").Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl11"
$ws.Cells.Item(2, 3).Value = "Ackr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.57077
$ws.Cells.Item(2, 8).Value = 4.71231
$ws.Cells.Item(2, 9).Value = 0.02582502173444737
$ws.Cells.Item(2, 10).Value = 0.02582502173444737
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.377371
$ws.Cells.Item(2, 14).Value = 1.132113
$ws.Cells.Item(2, 15).Value = 0.4698794580655765
$ws.Cells.Item(2, 16).Value = 0.4698794580655764
$ws.Cells.Item(2, 17).Value = 0.59276304567
$ws.Cells.Item(2, 18).Value = 5.33486741103
$ws.Cells.Item(2, 19).Value = 0.01213464721711386
$ws.Cells.Item(2, 20).Value = 0.01213464721711386

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl11"
$ws.Cells.Item(3, 3).Value = "Ackr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.57077
$ws.Cells.Item(3, 8).Value = 4.71231
$ws.Cells.Item(3, 9).Value = 0.02582502173444737
$ws.Cells.Item(3, 10).Value = 0.02582502173444737
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.3560133333333333
$ws.Cells.Item(3, 14).Value = 1.06804
$ws.Cells.Item(3, 15).Value = 0.443286188209444
$ws.Cells.Item(3, 16).Value = 0.443286188209444
$ws.Cells.Item(3, 17).Value = 0.5592150636000001
$ws.Cells.Item(3, 18).Value = 5.032935572400001
$ws.Cells.Item(3, 19).Value = 0.01144787544508922
$ws.Cells.Item(3, 20).Value = 0.01144787544508922

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ccl11"
$ws.Cells.Item(4, 3).Value = "Ackr4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.57077
$ws.Cells.Item(4, 8).Value = 4.71231
$ws.Cells.Item(4, 9).Value = 0.02582502173444737
$ws.Cells.Item(4, 10).Value = 0.02582502173444737
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.06973866666666667
$ws.Cells.Item(4, 14).Value = 0.209216
$ws.Cells.Item(4, 15).Value = 0.08683435372497944
$ws.Cells.Item(4, 16).Value = 0.08683435372497944
$ws.Cells.Item(4, 17).Value = 0.10954340544
$ws.Cells.Item(4, 18).Value = 0.9858906489600001
$ws.Cells.Item(4, 19).Value = 0.002242499072244285
$ws.Cells.Item(4, 20).Value = 0.002242499072244285

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ccl11"
$ws.Cells.Item(5, 3).Value = "Ackr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 56.85979966666667
$ws.Cells.Item(5, 8).Value = 170.579399
$ws.Cells.Item(5, 9).Value = 0.934831682683009
$ws.Cells.Item(5, 10).Value = 0.934831682683009
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.377371
$ws.Cells.Item(5, 14).Value = 1.132113
$ws.Cells.Item(5, 15).Value = 0.4698794580655765
$ws.Cells.Item(5, 16).Value = 0.4698794580655764
$ws.Cells.Item(5, 17).Value = 21.45723946000967
$ws.Cells.Item(5, 18).Value = 193.115155140087
$ws.Cells.Item(5, 19).Value = 0.4392582044416232
$ws.Cells.Item(5, 20).Value = 0.4392582044416232

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ccl11"
$ws.Cells.Item(6, 3).Value = "Ackr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 56.85979966666667
$ws.Cells.Item(6, 8).Value = 170.579399
$ws.Cells.Item(6, 9).Value = 0.934831682683009
$ws.Cells.Item(6, 10).Value = 0.934831682683009
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.3560133333333333
$ws.Cells.Item(6, 14).Value = 1.06804
$ws.Cells.Item(6, 15).Value = 0.443286188209444
$ws.Cells.Item(6, 16).Value = 0.443286188209444
$ws.Cells.Item(6, 17).Value = 20.24284681199556
$ws.Cells.Item(6, 18).Value = 182.18562130796
$ws.Cells.Item(6, 19).Value = 0.4143979732339716
$ws.Cells.Item(6, 20).Value = 0.4143979732339716

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ccl11"
$ws.Cells.Item(7, 3).Value = "Ackr4"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 56.85979966666667
$ws.Cells.Item(7, 8).Value = 170.579399
$ws.Cells.Item(7, 9).Value = 0.934831682683009
$ws.Cells.Item(7, 10).Value = 0.934831682683009
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.06973866666666667
$ws.Cells.Item(7, 14).Value = 0.209216
$ws.Cells.Item(7, 15).Value = 0.08683435372497944
$ws.Cells.Item(7, 16).Value = 0.08683435372497944
$ws.Cells.Item(7, 17).Value = 3.965326615687112
$ws.Cells.Item(7, 18).Value = 35.68793954118401
$ws.Cells.Item(7, 19).Value = 0.08117550500741413
$ws.Cells.Item(7, 20).Value = 0.08117550500741413

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Ccl11"
$ws.Cells.Item(8, 3).Value = "Ackr4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.9273763333333335
$ws.Cells.Item(8, 8).Value = 2.782129
$ws.Cells.Item(8, 9).Value = 0.01524698967025436
$ws.Cells.Item(8, 10).Value = 0.01524698967025436
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.377371
$ws.Cells.Item(8, 14).Value = 1.132113
$ws.Cells.Item(8, 15).Value = 0.4698794580655765
$ws.Cells.Item(8, 16).Value = 0.4698794580655764
$ws.Cells.Item(8, 17).Value = 0.3499649342863334
$ws.Cells.Item(8, 18).Value = 3.149684408577
$ws.Cells.Item(8, 19).Value = 0.007164247243390562
$ws.Cells.Item(8, 20).Value = 0.007164247243390561

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Ccl11"
$ws.Cells.Item(9, 3).Value = "Ackr4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.9273763333333335
$ws.Cells.Item(9, 8).Value = 2.782129
$ws.Cells.Item(9, 9).Value = 0.01524698967025436
$ws.Cells.Item(9, 10).Value = 0.01524698967025436
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.3560133333333333
$ws.Cells.Item(9, 14).Value = 1.06804
$ws.Cells.Item(9, 15).Value = 0.443286188209444
$ws.Cells.Item(9, 16).Value = 0.443286188209444
$ws.Cells.Item(9, 17).Value = 0.3301583396844445
$ws.Cells.Item(9, 18).Value = 2.971425057160001
$ws.Cells.Item(9, 19).Value = 0.006758779932595824
$ws.Cells.Item(9, 20).Value = 0.006758779932595824

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Ccl11"
$ws.Cells.Item(10, 3).Value = "Ackr4"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.9273763333333335
$ws.Cells.Item(10, 8).Value = 2.782129
$ws.Cells.Item(10, 9).Value = 0.01524698967025436
$ws.Cells.Item(10, 10).Value = 0.01524698967025436
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.06973866666666667
$ws.Cells.Item(10, 14).Value = 0.209216
$ws.Cells.Item(10, 15).Value = 0.08683435372497944
$ws.Cells.Item(10, 16).Value = 0.08683435372497944
$ws.Cells.Item(10, 17).Value = 0.06467398898488891
$ws.Cells.Item(10, 18).Value = 0.5820659008640001
$ws.Cells.Item(10, 19).Value = 0.001323962494267975
$ws.Cells.Item(10, 20).Value = 0.001323962494267975

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ccl11"
$ws.Cells.Item(11, 3).Value = "Ackr4"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.7810079999999999
$ws.Cells.Item(11, 8).Value = 2.343024
$ws.Cells.Item(11, 9).Value = 0.0128405486320577
$ws.Cells.Item(11, 10).Value = 0.0128405486320577
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.377371
$ws.Cells.Item(11, 14).Value = 1.132113
$ws.Cells.Item(11, 15).Value = 0.4698794580655765
$ws.Cells.Item(11, 16).Value = 0.4698794580655764
$ws.Cells.Item(11, 17).Value = 0.2947297699679999
$ws.Cells.Item(11, 18).Value = 2.652567929712
$ws.Cells.Item(11, 19).Value = 0.006033510032495949
$ws.Cells.Item(11, 20).Value = 0.006033510032495948

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ccl11"
$ws.Cells.Item(12, 3).Value = "Ackr4"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.7810079999999999
$ws.Cells.Item(12, 8).Value = 2.343024
$ws.Cells.Item(12, 9).Value = 0.0128405486320577
$ws.Cells.Item(12, 10).Value = 0.0128405486320577
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.3560133333333333
$ws.Cells.Item(12, 14).Value = 1.06804
$ws.Cells.Item(12, 15).Value = 0.443286188209444
$ws.Cells.Item(12, 16).Value = 0.443286188209444
$ws.Cells.Item(12, 17).Value = 0.27804926144
$ws.Cells.Item(12, 18).Value = 2.50244335296
$ws.Cells.Item(12, 19).Value = 0.005692037857622847
$ws.Cells.Item(12, 20).Value = 0.005692037857622847

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ccl11"
$ws.Cells.Item(13, 3).Value = "Ackr4"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.7810079999999999
$ws.Cells.Item(13, 8).Value = 2.343024
$ws.Cells.Item(13, 9).Value = 0.0128405486320577
$ws.Cells.Item(13, 10).Value = 0.0128405486320577
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.06973866666666667
$ws.Cells.Item(13, 14).Value = 0.209216
$ws.Cells.Item(13, 15).Value = 0.08683435372497944
$ws.Cells.Item(13, 16).Value = 0.08683435372497944
$ws.Cells.Item(13, 17).Value = 0.054466456576
$ws.Cells.Item(13, 18).Value = 0.490198109184
$ws.Cells.Item(13, 19).Value = 0.001115000741938899
$ws.Cells.Item(13, 20).Value = 0.001115000741938899

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Ccl11"
$ws.Cells.Item(14, 3).Value = "Ackr4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.6846153333333334
$ws.Cells.Item(14, 8).Value = 2.053846
$ws.Cells.Item(14, 9).Value = 0.01125575728023152
$ws.Cells.Item(14, 10).Value = 0.01125575728023152
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.377371
$ws.Cells.Item(14, 14).Value = 1.132113
$ws.Cells.Item(14, 15).Value = 0.4698794580655765
$ws.Cells.Item(14, 16).Value = 0.4698794580655764
$ws.Cells.Item(14, 17).Value = 0.2583539729553333
$ws.Cells.Item(14, 18).Value = 2.325185756598
$ws.Cells.Item(14, 19).Value = 0.005288849130952853
$ws.Cells.Item(14, 20).Value = 0.005288849130952852

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Ccl11"
$ws.Cells.Item(15, 3).Value = "Ackr4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.6846153333333334
$ws.Cells.Item(15, 8).Value = 2.053846
$ws.Cells.Item(15, 9).Value = 0.01125575728023152
$ws.Cells.Item(15, 10).Value = 0.01125575728023152
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.3560133333333333
$ws.Cells.Item(15, 14).Value = 1.06804
$ws.Cells.Item(15, 15).Value = 0.443286188209444
$ws.Cells.Item(15, 16).Value = 0.443286188209444
$ws.Cells.Item(15, 17).Value = 0.2437321868711111
$ws.Cells.Item(15, 18).Value = 2.19358968184
$ws.Cells.Item(15, 19).Value = 0.004989521740164529
$ws.Cells.Item(15, 20).Value = 0.004989521740164529

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Ccl11"
$ws.Cells.Item(16, 3).Value = "Ackr4"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.6846153333333334
$ws.Cells.Item(16, 8).Value = 2.053846
$ws.Cells.Item(16, 9).Value = 0.01125575728023152
$ws.Cells.Item(16, 10).Value = 0.01125575728023152
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.06973866666666667
$ws.Cells.Item(16, 14).Value = 0.209216
$ws.Cells.Item(16, 15).Value = 0.08683435372497944
$ws.Cells.Item(16, 16).Value = 0.08683435372497944
$ws.Cells.Item(16, 17).Value = 0.04774416052622223
$ws.Cells.Item(16, 18).Value = 0.429697444736
$ws.Cells.Item(16, 19).Value = 0.0009773864091141362
$ws.Cells.Item(16, 20).Value = 0.0009773864091141362

